$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1197518333792686, 0.9619067311286926, 0.01647702418267727, 0.9977655410766602)
    3 = @(0.0276475828140974, 0.9950518608093262, 0.00886818952858448, 0.9982123970985413)
    4 = @(0.01758970879018307, 0.9953539967536926, 0.002315406687557697, 0.9992551803588867)
    5 = @(0.01310524716973305, 0.9957506060600281, 0.00124630075879395, 0.9995530843734741)
    6 = @(0.01323879696428776, 0.9955050945281982, 0.001025006640702486, 0.9995530843734741)
    7 = @(0.01138939335942268, 0.9960905909538269, 0.0004758846189361066, 1)
    8 = @(0.0118631711229682, 0.9958261847496033, 0.000398884090827778, 1)
    9 = @(0.01131247542798519, 0.9961850047111511, 0.0004076385230291635, 1)
    10 = @(0.0100646149367094, 0.9966005086898804, 0.0002068919275188819, 1)
    11 = @(0.01094847824424505, 0.9961850047111511, 0.0001732196396915242, 1)
    12 = @(0.01059538591653109, 0.9961472153663635, 0.0003776536323130131, 0.9997020959854126)
    13 = @(0.01070791110396385, 0.9961661100387573, 0.0001256289688171819, 1)
    14 = @(0.01059665251523256, 0.9960716962814331, 0.0001186915542348288, 1)
    15 = @(0.01083364617079496, 0.9959017038345337, 0.00007473116420442238, 1)
    16 = @(0.01033419277518988, 0.9960528016090393, 0.0002333047305000946, 1)
    17 = @(0.01053417753428221, 0.9959583878517151, 0.00008618435822427273, 1)
    18 = @(0.01055614184588194, 0.9960528016090393, 0.0001234213705174625, 1)
    19 = @(0.009628799743950367, 0.996543824672699, 0.0001448042748961598, 1)
    20 = @(0.009384415112435818, 0.9965060949325562, 0.0001065137403202243, 1)
    21 = @(0.01095025427639484, 0.9961850047111511, 0.0000796353560872376, 1)
    22 = @(0.00979688111692667, 0.996430516242981, 0.0000154764511535177, 1)
    23 = @(0.01047691982239485, 0.9961094856262207, 0.00006192500586621463, 1)
    24 = @(0.01068889629095793, 0.9960150122642517, 0.00006898287392687052, 1)
    25 = @(0.009886534884572029, 0.9963927268981934, 0.00006761521217413247, 1)
    26 = @(0.01028876006603241, 0.9963927268981934, 0.00004689432171289809, 1)
    27 = @(0.01087967865169048, 0.9959017038345337, 0.00001347345460089855, 1)
    28 = @(0.01002767868340015, 0.9960716962814331, 0.0000562663481105119, 1)
    29 = @(0.009760917164385319, 0.9962227940559387, 0.000006476855105574941, 1)
    30 = @(0.01034883037209511, 0.9961094856262207, 0.00004765309495269321, 1)
    31 = @(0.01053190790116787, 0.9959017038345337, 0.00004065906614414416, 1)
    32 = @(0.009352984838187695, 0.996638298034668, 0.00001336712739430368, 1)
    33 = @(0.0106295133009553, 0.9959583878517151, 0.00003139493855996989, 1)
    34 = @(0.008538245223462582, 0.99709153175354, 0.00001219619753101142, 1)
    35 = @(0.01035260781645775, 0.9961094856262207, 0.000009434774256078526, 1)
    36 = @(0.01062451489269733, 0.9961472153663635, 0.00002635655800986569, 1)
    37 = @(0.00994833093136549, 0.9963549971580505, 0.000003242114644308458, 1)
    38 = @(0.009668433107435703, 0.9964494109153748, 0.000006651015610259492, 1)
    39 = @(0.009725398384034634, 0.9962794184684753, 0.000009608293112250976, 1)
    40 = @(0.009613302536308765, 0.9964116215705872, 0.000005880316621187376, 1)
    41 = @(0.009768443182110786, 0.9963549971580505, 0.000006880087312310934, 1)
    42 = @(0.01015400514006615, 0.9963172078132629, 0.000007980209375091363, 1)
    43 = @(0.01133520063012838, 0.9959394931793213, 0.00001384930783387972, 1)
    44 = @(0.009894312359392643, 0.9962983131408691, 0.00001271406199521152, 1)
    45 = @(0.01062796264886856, 0.9958828091621399, 0.00003416154504520819, 1)
    46 = @(0.00947289913892746, 0.9963927268981934, 0.00001839303331507836, 1)
    47 = @(0.009208904579281807, 0.996657133102417, 0.00002657950062712189, 1)
    48 = @(0.009058337658643723, 0.9966005086898804, 0.00002032283919106703, 1)
    49 = @(0.009793973527848721, 0.9964683055877686, 0.00001378701381327119, 1)
    50 = @(0.01044525112956762, 0.9960528016090393, 0.00001425484515493736, 1)
    51 = @(0.009618457406759262, 0.9966194033622742, 0.00002051459887297824, 1)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Cells.Item([int]$r, $c + 1).Value = $rowVals[$c]
    }
}
